# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet (fund holdings detail) positioned right
#    before the "总计" (totals) sheet.
# 2. Insert a new top row into "总计" summarising the 2022-Q1 quarter and
#    renumber the existing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q1" worksheet
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")

$new = $wb.Worksheets.Add()
$new.Name = "2022-Q1"

# Carry over the header-row / index-column look (style only) from the
# "2021-Q4" sheet so the new sheet matches the existing quarterly sheets.
# (Column A has no header cell in row 1, matching the other quarterly sheets.)
$src.Range("B1:H1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)
$src.Range("A2").Copy()
$new.Range("A2:A11").PasteSpecial(-4122)

# Header row
$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

# Fund holding rows. Columns B (code) and D/E/F/G (numeric-looking figures
# taken verbatim from the source report) are kept as text -- formatting the
# cell as Text first stops Excel re-interpreting "005888" / "0.9638" etc as
# numbers (which would also strip the leading zeros on the fund codes).
$rows = @(
    @{A=0; B="160311"; C="华夏蓝筹混合(LOF)";            D="30.12"; E="87.73"; F="3.20"; G="0.9638"; H=9}
    @{A=1; B="005888"; C="华夏新兴消费混合A";              D="16.37"; E="91.96"; F="5.43"; G="0.8889"; H=1}
    @{A=2; B="012421"; C="华夏优加生活混合A";              D="8.67";  E="92.98"; F="5.78"; G="0.5011"; H=1}
    @{A=3; B="001479"; C="中邮风格轮动灵活配置混合";        D="9.45";  E="62.17"; F="3.59"; G="0.3393"; H=4}
    @{A=4; B="005889"; C="华夏新兴消费混合C";              D="4.25";  E="91.96"; F="5.43"; G="0.2308"; H=1}
    @{A=5; B="005457"; C="景顺长城量化小盘股票";            D="9.49";  E="93.39"; F="1.96"; G="0.1860"; H=3}
    @{A=6; B="012422"; C="华夏优加生活混合C";              D="0.17";  E="92.98"; F="5.78"; G="0.0098"; H=1}
    @{A=7; B="004641"; C="万家量化睿选灵活配置混合";        D="0.16";  E="85.90"; F="1.36"; G="0.0022"; H=5}
    @{A=8; B="008437"; C="九泰行业优选灵活配置混合A";       D="0.11";  E="51.13"; F="1.38"; G="0.0015"; H=10}
    @{A=9; B="008438"; C="九泰行业优选灵活配置混合C";       D="0.06";  E="51.13"; F="1.38"; G="0.0008"; H=10}
)

$r = 2
foreach ($row in $rows) {
    $new.Range("A$r").Value = $row.A

    $cB = $new.Range("B$r")
    $cB.NumberFormat = "@"
    $cB.Value = $row.B

    $new.Range("C$r").Value = $row.C

    $cD = $new.Range("D$r")
    $cD.NumberFormat = "@"
    $cD.Value = $row.D

    $cE = $new.Range("E$r")
    $cE.NumberFormat = "@"
    $cE.Value = $row.E

    $cF = $new.Range("F$r")
    $cF.NumberFormat = "@"
    $cF.Value = $row.F

    $cG = $new.Range("G$r")
    $cG.NumberFormat = "@"
    $cG.Value = $row.G

    $new.Range("H$r").Value = $row.H

    $r = $r + 1
}

# Move the freshly populated sheet into place, right before "总计".
$new.Move($null, $wb.Worksheets.Item("2021-Q4"))

# ---------------------------------------------------------------------
# 2) Prepend the 2022-Q1 summary row to "总计"
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

$tot.Rows.Item(2).Insert()
$tot.Rows.Item(2).ClearFormats()

# Re-use the existing index-column style (bold/centred/bordered) instead of
# whatever formatting Insert() guessed for the new row.
$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 10
$tot.Range("D2").Value = 3.12

# Renumber the index column for the rows that shifted down.
$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
$tot.Range("A6").Value = 4
$tot.Range("A7").Value = 5
